# Slide 3 contains the security checklist table (shape 2, "Group 518").
# Row 4 / Column 4 holds the "X" (적용여부) marker for the password
# encryption requirement; change it to "O" since the password is now
# stored as a SHA256 hash. Row 4 / Column 5 (비고 / remarks) gets the
# note "SHA256" describing the new hashing scheme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tableShape = $s.Shapes.Item(2)
$tbl = $tableShape.Table

$statusCell = $tbl.Cell(4, 4)
$statusCell.Shape.TextFrame.TextRange.Text = "O"

$remarkCell = $tbl.Cell(4, 5)
$remarkCell.Shape.TextFrame.TextRange.Text = "SHA256"
